$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before Q (which shifts old Q,R,S -> R,S,T)
$ws.Range("Q1").EntireColumn.Insert()

# Insert a new column before S (after the first insert, old R (peso_vab_enge) now sits
# at S; we need the new "brecha_productividad_sin_enge_tot_level" column before it)
$ws.Range("S1").EntireColumn.Insert()

# Header row
$ws.Range("Q1").Value = "productividad_sin_enge_level"
$ws.Range("S1").Value = "brecha_productividad_sin_enge_tot_level"

# Data for new Q column (productividad_sin_enge_level)
$qValues = @{
  2 = 0.03106842400804815
  3 = 0.03446543551152889
  4 = 0.03549292869425207
  5 = 0.03957096078447744
  6 = 0.04091524678415288
  7 = 0.03822206121711909
  8 = 0.03723843319785487
  9 = 0.03419106418387535
  10 = 0.03091150546538797
  11 = 0.02015463780085779
  12 = 0.02290597434741486
  13 = 0.02326513082716702
  14 = 0.02275982697302471
  15 = 0.02396416960311077
  16 = 0.02434398946512564
  17 = 0.02358195762308805
  18 = 0.02321566860077521
  19 = 0.0295672963153145
  20 = 0.03007716557868893
  21 = 0.02582327438604886
  22 = 0.02626994052581537
  23 = 0.02321446006483791
  24 = 0.02411999307866832
  25 = 0.02395886545443953
  26 = 0.02010442393166966
  27 = 0.01804310344368463
  28 = 0.01752299618113438
  29 = 0.02684010311130477
}

# Data for new S column (brecha_productividad_sin_enge_tot_level)
$sValues = @{
  2 = 0.7217866279342013
  3 = 0.7061347470885972
  4 = 0.7205939142286747
  5 = 0.7274295185965005
  6 = 0.7158631339776839
  7 = 0.6964979768282287
  8 = 0.7044691905273722
  9 = 0.6572341654741737
  10 = 0.640122874364606
  11 = 0.4664758337295752
  12 = 0.5188359339706339
  13 = 0.4876241351803389
  14 = 0.4848687707668133
  15 = 0.4933955560261969
  16 = 0.4607912212026656
  17 = 0.4512128566514683
  18 = 0.4811258478933602
  19 = 0.5071842261070042
  20 = 0.5093578642283374
  21 = 0.473622428045198
  22 = 0.4723300233404075
  23 = 0.4421653496897042
  24 = 0.4454741915504722
  25 = 0.4774814540996213
  26 = 0.4121049083080825
  27 = 0.3820579528613384
  28 = 0.4268772236080584
  29 = 0.5274793916539007
}

foreach ($row in 2..29) {
    $ws.Range("Q$row").Value = $qValues[$row]
    $ws.Range("S$row").Value = $sValues[$row]
}
